$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.411.35"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.848.02"
$ws.Range("E3").Value = "  -0.14%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.02%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.70"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -1.05%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6321"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -3.90%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.346.01"
$ws.Range("E8").Value = "  +80.97%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07572"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +1.13%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2968"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -1.20%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.58"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "3.510.06"
$ws.Range("E12").Value = "  +68.45%  "
$ws.Range("E13").Value = "  +0.78%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.982"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -0.84%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6835"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.30%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009972"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +4.89%  "
$ws.Range("E17").Value = "  -1.10%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.161"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "29.441.95"
$ws.Range("E19").Value = "  -0.17%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.29"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -2.61%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.49"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -0.61%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +0.03%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.564"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("E24").Value = "  +0.01%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.28"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -0.97%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1384"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -3.01%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.427"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.77%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.65"
$ws.Range("D28").Style = $style
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.468"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.38%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05798"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -2.99%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.257"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.20%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.124"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").Value = "3.505.16"
$ws.Range("E34").Value = "  +75.29%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.866"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -1.77%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7161"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -1.01%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.593"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "1.247.88"
$ws.Range("E39").Value = "  +3.81%  "
$ws.Range("E40").Value = "  -0.10%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01805"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +1.45%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9018"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -0.94%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.099"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  +0.04%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.56"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -0.29%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.83"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.48%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.194"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -3.36%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.164"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +0.59%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4011"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.07%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.686"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +1.24%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1123"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -0.45%  "
